$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 347 - this shifts the existing rows 347-372
# (and all their data) down to rows 348-373, preserving them unchanged.
$ws.Rows.Item(347).Insert()

# Populate the newly inserted row 347 with the new weekly price record.
$ws.Cells.Item(347, 1).Value = 10
$ws.Cells.Item(347, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(347, 3).Value = "La Araucanía"
$ws.Cells.Item(347, 4).Value = 44826
$ws.Cells.Item(347, 5).Value = 9
$ws.Cells.Item(347, 6).Value = 100112017
$ws.Cells.Item(347, 7).Value = "Apio"
$ws.Cells.Item(347, 8).Value = "Americana (o)"
$ws.Cells.Item(347, 9).Value = "Primera"
$ws.Cells.Item(347, 10).Value = 200
$ws.Cells.Item(347, 11).Value = 12000
$ws.Cells.Item(347, 12).Value = 12000
$ws.Cells.Item(347, 13).Value = 12000
$ws.Cells.Item(347, 14).Value = "$/docena de matas"
$ws.Cells.Item(347, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(347, 16).Value = 2000
$ws.Cells.Item(347, 17).Value = 6
$ws.Cells.Item(347, 18).Value = "Hortaliza"
